$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Remove the old "last updated" date stamp in C1.
$about.Range("C1").Clear()

# Insert 4 new rows for the updated "Sources:" citation block (year,
# document title, URL, table reference) right under the "Sources:" row.
$about.Rows("4:7").Insert()
$about.Range("B3").Value = "United States EPA"
$about.Range("B4").Value = 2012
$about.Range("B4").HorizontalAlignment = -4131
$about.Range("B5").Value = "Consumer Vehicle Choice Model Documentation"
$about.Range("B6").Value = "https://nepis.epa.gov/Exe/ZyPDF.cgi/P100EZ37.PDF?Dockey=P100EZ37.PDF"
$about.Range("B7").Value = "Table 5 Generalized Cost Coefficient Calibration"

# Insert 3 new rows for the extra explanatory note about the -3 / -5
# calibration choice, right after the existing "Notes" paragraph.
$about.Rows("16:18").Insert()
$about.Range("A16").Value = "We choose a value of -3 for passenger vehicles and a value of -5 for other vehicle types, "
$about.Range("A17").Value = "based on the ranges in Table 5 of the cited EPA documentation."

# ---------------------------------------------------------------------
# "TTLE" sheet - update logit exponents from -3 to -5
# ---------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")
$ttle.Range("B2:C7").Value = -5
